# Applies the "Used sonarlint, fixed some problems" commit to the workbook.
$wb = $excel.ActiveWorkbook

# --- View-state clean-up on the first three sheets: Excel had them scrolled away
#     from A1 and/or marked as the tab shown on open; reset that. ---
$wsReq = $wb.Worksheets.Item("Requirements Phase Defects")
$wsReq.Range("H16").Select()

$wsArch = $wb.Worksheets.Item("Architect. Design Phase Defects")
$wsArch.Range("E28").Select()

$wsCode = $wb.Worksheets.Item("Coding Phase Defects")
$wsCode.Range("E32").Select()

# --- Tool-basedCodeAnalysis sheet: fill in the SonarLint findings table ---
$ws = $wb.Worksheets.Item("Tool-basedCodeAnalysis")

# Tool used: SonarLint
$ws.Range("D4").Value = "SonarLint"

# Row 10 - clone() override finding (filled D, E, C, F order)
$ws.Range("D10").Value = "clone should not be overridden"
$ws.Range("E10").Value = "clone method was overridden"
$ws.Range("C10").Value = "Task.java, 175"
$ws.Range("F10").Value = "clone method override implementation was removed"
$ws.Range("F10").Font.Italic = $true
$ws.Range("F10").Borders.LineStyle = -4142
$ws.Range("F10").WrapText = $false

# Row 11 - constant naming convention finding
$ws.Range("C11").Value = "TaskIO.java, 18, 19, 20"
$ws.Range("D11").Value = "Constant names should comply with a naming convention"
$ws.Range("E11").Value = "secondsInDay, secondsInHour, secondsInMin"
$ws.Range("F11").Value = "SECONDS_IN_DAY, SECONDS_IN_HOUR, SECONDS_IN_MIN"

# Row 12 - try-with-resources finding
$ws.Range("C12").Value = "TaskIO.java, 25, 47, 73, 87, 123, 136"
$ws.Range("D12").Value = "Try-with-resources should be used"
$ws.Range("E12").Value = "finally clause used"
$ws.Range("F12").Value = "used try-with-resources + solved possible NullPointerExceptions"

# Row 13 - utility class constructor finding
$ws.Range("C13").Value = "TaskIO.java, 15"
$ws.Range("D13").Value = "Utility classes should not have public constructors"
$ws.Range("E13").Value = "no explicit constructor"
$ws.Range("F13").Value = "added private constructor"

# Row heights to fit the wrapped comments
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 90
$ws.Rows.Item(12).RowHeight = 75
$ws.Rows.Item(13).RowHeight = 45

# Selection / view state: this sheet becomes the active one
$ws.Range("G15").Select()
$excel.ActiveWindow.Zoom = 115
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1

$ws.Activate()
